# Revised estimates for the updated paper draft.
# The table gains a new column (E) and every numeric estimate in the
# body of the table is refreshed; the old text header labels in row 1
# are replaced by the (now unlabeled) numeric placeholders that the
# authors' export script produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row) ------------------------------------------------
# Former header cells (shared-string labels) become plain numbers.
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 0

# New column E, header cell, reuse the same formatting (thin border,
# bold, centered/top-aligned) already used by the rest of row 1.
$ws.Range("E1").Value = 0
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 2: FE -----------------------------------------------------------
$ws.Range("A2").Value = "FE"
$ws.Range("B2").Value = 1.38
$ws.Range("C2").Value = 0.24
$ws.Range("D2").Value = 0.2

# --- Row 3: FE+Disg -------------------------------------------------------
$ws.Range("A3").Value = "FE+Disg"
$ws.Range("B3").Value = 0.71
$ws.Range("C3").Value = 0.24
$ws.Range("D3").Value = 0.2

# --- Row 4: FE+Disg+Var ----------------------------------------------------
$ws.Range("A4").Value = "FE+Disg+Var"
$ws.Range("B4").Value = 0.95
$ws.Range("C4").Value = 0.21
$ws.Range("D4").Value = 0.2

$wb.Save()
